# "xml file undate in pom"
# Sheet3!D4 held an order/tracking-style numeric-looking text value
# (shared string "000001515047221824"); update it to the new value
# "000001515240526631", keeping it stored as text (not auto-coerced to a
# number, which would also strip the leading zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$cell = $ws.Range("D4")

# Force text interpretation so the long numeric-looking string keeps its
# leading zeros and is written back as a shared string (t="s"), not a
# numeric cell. Restore the number format afterwards so no new cell style
# is introduced (the source cell carries no explicit style either).
$cell.NumberFormat = "@"
$cell.Value = "000001515240526631"
$cell.NumberFormat = "General"
